$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-46 (age 20-64): B,C,D columns
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 2).Value = 0.000417
    $ws.Cells.Item($r, 3).Value = 0.000805
    $ws.Cells.Item($r, 4).Value = 0.001527
}

# Update rows 47-87 (age 65-105): B,C,D columns
for ($r = 47; $r -le 87; $r++) {
    $ws.Cells.Item($r, 2).Value = 0.005813
    $ws.Cells.Item($r, 3).Value = 0.011223
    $ws.Cells.Item($r, 4).Value = 0.021296
}

# Update the active cell selection
$ws.Range("J11").Select()
